$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 19179
$ws1.Range("F5").Value = 144
$ws1.Range("F15").Value = 232
$ws1.Range("F22").Value = 8007
$ws1.Range("F36").Value = 820

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 19179
$ws4.Range("F5").Value = 144
$ws4.Range("F15").Value = 232
$ws4.Range("F23").Value = 8007
$ws4.Range("F39").Value = 820
